# Insert a new weekly data point for "Feria Lagunitas de Puerto Montt - Perejil".
# The sheet is ordered from most-recent date (row 164) to oldest (row 175).
# A new, more-recent record is inserted at row 164, and all the previously
# existing rows 164-175 shift down by one to become rows 165-176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that stay constant across every data row of this sheet.
$colA = 4
$colB = "Feria Lagunitas de Puerto Montt"
$colC = "Los Lagos"
$colE = 10
$colF = 100112044
$colG = "Perejil"
$colH = "Sin especificar"
$colI = "Primera"
$colR = "Hortaliza"

# Work from the bottom up so we never overwrite a source row before reading it.
for ($r = 175; $r -ge 164; $r--) {
    $target = $r + 1

    $d = $ws.Cells.Item($r, 4).Value2
    $j = $ws.Cells.Item($r, 10).Value2
    $k = $ws.Cells.Item($r, 11).Value2
    $l = $ws.Cells.Item($r, 12).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    $n = $ws.Cells.Item($r, 14).Value2
    $o = $ws.Cells.Item($r, 15).Value2
    $p = $ws.Cells.Item($r, 16).Value2
    $q = $ws.Cells.Item($r, 17).Value2

    $ws.Cells.Item($target, 1).Value2 = $colA
    $ws.Cells.Item($target, 2).Value2 = $colB
    $ws.Cells.Item($target, 3).Value2 = $colC
    $ws.Cells.Item($target, 4).Value2 = $d
    $ws.Range("D" + $target).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($target, 5).Value2 = $colE
    $ws.Cells.Item($target, 6).Value2 = $colF
    $ws.Cells.Item($target, 7).Value2 = $colG
    $ws.Cells.Item($target, 8).Value2 = $colH
    $ws.Cells.Item($target, 9).Value2 = $colI
    $ws.Cells.Item($target, 10).Value2 = $j
    $ws.Cells.Item($target, 11).Value2 = $k
    $ws.Cells.Item($target, 12).Value2 = $l
    $ws.Cells.Item($target, 13).Value2 = $m
    $ws.Cells.Item($target, 14).Value2 = $n
    $ws.Cells.Item($target, 15).Value2 = $o
    $ws.Cells.Item($target, 16).Value2 = $p
    $ws.Cells.Item($target, 17).Value2 = $q
    $ws.Cells.Item($target, 18).Value2 = $colR
}

# New record inserted at row 164.
$ws.Cells.Item(164, 1).Value2 = $colA
$ws.Cells.Item(164, 2).Value2 = $colB
$ws.Cells.Item(164, 3).Value2 = $colC
$ws.Cells.Item(164, 4).Value2 = 44516
$ws.Range("D164").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(164, 5).Value2 = $colE
$ws.Cells.Item(164, 6).Value2 = $colF
$ws.Cells.Item(164, 7).Value2 = $colG
$ws.Cells.Item(164, 8).Value2 = $colH
$ws.Cells.Item(164, 9).Value2 = $colI
$ws.Cells.Item(164, 10).Value2 = 180
$ws.Cells.Item(164, 11).Value2 = 5000
$ws.Cells.Item(164, 12).Value2 = 5000
$ws.Cells.Item(164, 13).Value2 = 5000
$ws.Cells.Item(164, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(164, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(164, 16).Value2 = 1667
$ws.Cells.Item(164, 17).Value2 = 3
$ws.Cells.Item(164, 18).Value2 = $colR
